# Update the cryptocurrency "Price" (column D) and "Volume(1h)" (column E)
# values on the active sheet to reflect the latest scrape, matching the
# GitHub Actions commit "Updated cryptos list ...".
#
# Note: several Price values look like plain decimal numbers (e.g. "0.7116").
# A leading apostrophe is used so Excel stores them as literal text (matching
# the original text-formatted cells) instead of auto-converting them to
# numbers, which would lose formatting like trailing zeros.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.384.81"
$ws.Range("E2").Value = "  +0.43%  "
$ws.Range("D3").Value = "1.876.34"
$ws.Range("E3").Value = "  +0.66%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "'0.7116"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").Value = "'242.04"
$ws.Range("E6").Value = "  +0.42%  "
$ws.Range("D7").Value = "'0.9999"
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("D8").Value = "'0.3115"
$ws.Range("E8").Value = "  +0.97%  "
$ws.Range("D9").Value = "'0.07791"
$ws.Range("E9").Value = "  +0.63%  "
$ws.Range("E10").Value = "  +0.89%  "
$ws.Range("D11").Value = "'0.08457"
$ws.Range("E11").Value = "  +1.36%  "
$ws.Range("D12").Value = "1.880.59"
$ws.Range("E12").Value = "  +0.19%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("D14").Value = "'0.7123"
$ws.Range("E14").Value = "  -0.42%  "
$ws.Range("D15").Value = "'91.22"
$ws.Range("E15").Value = "  +0.25%  "
$ws.Range("D16").Value = "29.387.24"
$ws.Range("E16").Value = "  +0.32%  "
$ws.Range("D17").Value = "'6.054"
$ws.Range("E17").Value = "  +1.50%  "
$ws.Range("D18").Value = "'0.000008237"
$ws.Range("E18").Value = "  +5.21%  "
$ws.Range("D19").Value = "'240.99"
$ws.Range("E19").Value = "  -0.82%  "
$ws.Range("E20").Value = "  +0.52%  "
$ws.Range("D21").Value = "2.118.16"
$ws.Range("E21").Value = "  -0.63%  "
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").Value = "'7.781"
$ws.Range("E23").Value = "  -1.64%  "
$ws.Range("D24").Value = "'0.9998"
$ws.Range("E24").Value = "  +0.01%  "
$ws.Range("D25").Value = "'0.1595"
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").Value = "'163.37"
$ws.Range("E26").Value = "  -0.06%  "
$ws.Range("E27").Value = "  +1.89%  "
$ws.Range("D28").Value = "'18.52"
$ws.Range("E28").Value = "  -0.32%  "
$ws.Range("E29").Value = "  +0.94%  "
$ws.Range("D30").Value = "'4.432"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").Value = "'4.324"
$ws.Range("E31").Value = "  +1.50%  "
$ws.Range("E32").Value = "  -4.48%  "
$ws.Range("D33").Value = "'0.05290"
$ws.Range("E33").Value = "  +2.53%  "
$ws.Range("D34").Value = "'1.943"
$ws.Range("E34").Value = "  +0.59%  "
$ws.Range("D35").Value = "'1.179"
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("D36").Value = "'0.7435"
$ws.Range("E36").Value = "  -11.97%  "
$ws.Range("D37").Value = "'2.696"
$ws.Range("E37").Value = "  +0.53%  "
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("D39").Value = "1.225.50"
$ws.Range("E39").Value = "  +4.94%  "
$ws.Range("D40").Value = "'2.725"
$ws.Range("E40").Value = "  +1.26%  "
$ws.Range("D41").Value = "'6.487"
$ws.Range("E41").Value = "  +4.62%  "
$ws.Range("D42").Value = "'0.8934"
$ws.Range("E42").Value = "  -0.27%  "
$ws.Range("D43").Value = "'110.08"
$ws.Range("E43").Value = "  +7.69%  "
$ws.Range("D44").Value = "'72.89"
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "'0.9998"
$ws.Range("D46").Value = "2.015.14"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("D47").Value = "'1.817"
$ws.Range("E47").Value = "  +1.57%  "
$ws.Range("D48").Value = "'0.5207"
$ws.Range("E48").Value = "  +0.50%  "
$ws.Range("E49").Value = "  +1.99%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").Value = "'0.4329"
$ws.Range("E51").Value = "  +1.12%  "
